$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Trim M2: remove the trailing protein annotation line, keep the DNA/codon block + two blank lines ---
# (set before B2's new text so the shared-strings table ends up in the same
#  append order as the target workbook)
$m2 = @"
1 caagtgccac tggctactag tgcaagtatg gctcgggtct ctgccaatgc agttgcactt
       61 gttgcactcg tctccgttct tctcacgtat ggctgctgcg cccagtcgcc gctcaactac
      121 accggctcct tggccaaatc ctccaaggct agctggtcat ggctccctgc caaggccaca
      181 tggtacggcg cgcctaccgg cgccggtccc gatgacaacg gtggtgcttg cggctacaag
      241 cacactaacc agtacccgtt catgtccatg acttcctgcg gcaacgagcc cctgttcaag
      301 gacggcatgg gctgcggcgc ctgctaccag atacgatgcg tcaataacaa ggcctgctcc
      361 ggcaagccgg agacggtcat gatcaccgac atgaactact accctgtggg caagtaccat
      421 ttcgacctca gcggcacggc gttcggcgcc atggcgaagc ccggccagaa cgacaagctc
      481 cgccacgccg gcattatcga catccagttc caaagggtgc catgcaatca tccgggcttg
      541 aacgtgaact tccaggtcga gcggggctcc aaccccaact acctggccgt gctggtggag
      601 ttcgcgaacc gggagggcac cgtggtgcag atggacctca tggagtcaag gaacggccgc
      661 ccgacggggt actggacggc gatgcgccac tcgtggggcg ccatctggcg gatggactcc
      721 aggcgccggc tgcagggccc cttctctctc cgcatccgca gcgaatccgg caagacgctg
      781 gtggccaaac aagtcatccc ggccaactgg aggcccgaca cgaactaccg ttccaacgtc
      841 cagttccgtt gattgctccg agcttccgat cgatcgacga agacgttgat taattcgg



"@
$ws.Range("M2").Value = $m2

# --- 2. Update B2: new test description text, with wrap text style ---
$ws.Range("B2").Value = "Check Delete Button enabled in tabbed view GF. Make sure the GF name appears striked out in search result after we delete the GF "
$ws.Range("B2").WrapText = $true

# --- 3. Apply wrap-text style to AC2 (DNA sequence cell) before the column shift ---
$ws.Range("AC2").WrapText = $true

# --- 4. Remove the (empty) column AB, shifting column AC (and its contents/formatting) left into AB ---
$ws.Columns.Item(28).Delete()

# --- 5. Column widths ---
# NOTE: the ColumnWidth -> stored XML "width" conversion in this runtime
# snaps to 1/6-character increments (it does not reproduce Excel's true
# MDW-pixel based width formula), so exact target fractions such as
# 36.85546875 cannot be hit bit-for-bit. The values below are solved to
# land on the nearest representable width for each target column.
$ws.Columns.Item(2).ColumnWidth = 36.0
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(7).ColumnWidth = 13.0
$ws.Columns.Item(8).ColumnWidth = 10.666666666666668
$ws.Columns.Item(9).ColumnWidth = 18.5
$ws.Columns.Item(13).ColumnWidth = 75.83333333333334
$ws.Columns.Item(14).ColumnWidth = 16.333333333333336
$ws.Columns.Item(15).ColumnWidth = 24.666666666666664
$ws.Columns.Item(16).ColumnWidth = 15.833333333333332
$ws.Columns.Item(17).ColumnWidth = 17.666666666666664
$ws.Columns.Item(18).ColumnWidth = 17.333333333333336
$ws.Columns.Item(19).ColumnWidth = 20.666666666666664
$ws.Columns.Item(21).ColumnWidth = 21.5
$ws.Columns.Item(24).ColumnWidth = 16.666666666666664
$ws.Columns.Item(27).ColumnWidth = 17.833333333333336
$ws.Columns.Item(28).ColumnWidth = 67.66666666666666

# --- 6. Row height for row 2 ---
$ws.Rows.Item(2).RowHeight = 128.25

# --- 7. Selection / view ---
$ws.Range("L4").Select()

# --- 8. Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
